# Rename default cover file: "cover.jpg" -> "default_cover.jpg" on the
# Metadata sheet (cover_path column), and leave the workbook with the
# Metadata sheet active/selected (as it was when the edit was made).

$wb = $excel.ActiveWorkbook
$wsWords = $wb.Worksheets.Item("Word entries")
$wsMeta  = $wb.Worksheets.Item("Metadata")

# Update the cover file name referenced in the metadata table.
$wsMeta.Range("E2").Value = "default_cover.jpg"

# Minor column-width touch-up on "Word entries" column D (re-saved by the
# authoring tool alongside the cover-file rename).
$wsWords.Columns.Item(4).ColumnWidth = 50.83

# The edit was made on the Metadata sheet, so it becomes the active tab.
$wsMeta.Activate()
$wsMeta.Range("E3").Select()
